# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" -> "Impact" bullet list so each
# line reads as an impact-focused accomplishment statement, and drop the
# final "Provided expert testimony..." bullet entirely (6 bullets -> 5).
#
# We anchor on the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph so the
# edits only touch that section and not the near-duplicate bullet text that
# also appears earlier under "Partner - Siege Analytics".

$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph.
$anchor = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("KEY ACHIEVEMENTS AND IMPACT")) {
        $anchor = $i
        break
    }
}
if ($anchor -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# Bullets live right after the "Impact" sub-heading: anchor+1 is "Impact",
# anchor+2 .. anchor+7 are the six bullet paragraphs.
$b1 = $anchor + 2
$b2 = $anchor + 3
$b3 = $anchor + 4
$b4 = $anchor + 5
$b5 = $anchor + 6
$b6 = $anchor + 7

# Replace bullet text in place (search scoped to each paragraph's own Range
# so the near-duplicate wording elsewhere in the document is untouched).
$d.Paragraphs.Item($b1).Range.Find.Execute(
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions", 2) | Out-Null

$d.Paragraphs.Item($b2).Range.Find.Execute(
    "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• 178% accuracy improvement in racial classification algorithms", 2) | Out-Null

$d.Paragraphs.Item($b3).Range.Find.Execute(
    "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%", 2) | Out-Null

$d.Paragraphs.Item($b4).Range.Find.Execute(
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• `$4.7M savings enabled nonprofit access", 2) | Out-Null

$d.Paragraphs.Item($b5).Range.Find.Execute(
    "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations", 2) | Out-Null

# Remove the trailing "Provided expert testimony..." bullet paragraph outright.
$d.Paragraphs.Item($b6).Range.Delete() | Out-Null

Write-Output "Done."
